# The commit swaps the two theme parts in the package:
#   - ppt/theme/theme1.xml ("Office Theme" / "Office" palette) <-> 
#     ppt/theme/theme2.xml ("Integral" / "Red Violet" palette)
# theme2.xml is the theme actually applied to the deck's (only) slide
# master, i.e. the one reachable/editable through the PowerPoint object
# model as SlideMaster.ColorScheme. Re-point each of its 12 scheme
# colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) at the
# "Office" palette values that theme1.xml carried before the edit, so
# the presentation's applied design becomes the plain Office theme.

function ColorRGB($r, $g, $b) {
    # PowerPoint's RGBColor.RGB uses the Win32 COLORREF packing
    # (0x00BBGGRR), i.e. blue in the high byte, red in the low byte.
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.ColorScheme

$colors.Colors(1).RGB  = ColorRGB 0x00 0x00 0x00   # dk1
$colors.Colors(2).RGB  = ColorRGB 0xFF 0xFF 0xFF   # lt1
$colors.Colors(3).RGB  = ColorRGB 0x44 0x54 0x6A   # dk2
$colors.Colors(4).RGB  = ColorRGB 0xE7 0xE6 0xE6   # lt2
$colors.Colors(5).RGB  = ColorRGB 0x5B 0x9B 0xD5   # accent1
$colors.Colors(6).RGB  = ColorRGB 0xED 0x7D 0x31   # accent2
$colors.Colors(7).RGB  = ColorRGB 0xA5 0xA5 0xA5   # accent3
$colors.Colors(8).RGB  = ColorRGB 0xFF 0xC0 0x00   # accent4
$colors.Colors(9).RGB  = ColorRGB 0x44 0x72 0xC4   # accent5
$colors.Colors(10).RGB = ColorRGB 0x70 0xAD 0x47   # accent6
$colors.Colors(11).RGB = ColorRGB 0x05 0x63 0xC1   # hlink
$colors.Colors(12).RGB = ColorRGB 0x95 0x4F 0x72   # folHlink
